# Update the crypto price list as produced by the "Updated symbol list"
# GitHub Actions workflow run on Fri Dec 16 22:33:39 UTC 2022.
#
# Column D holds prices that are stored as *text* (not numbers) in the
# workbook, so every numeric-looking value is written with a leading
# apostrophe to force Excel to keep it as text instead of silently
# converting it to a Double.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

# Simple price refreshes (column D only)
Set-TextValue "D2"  "241.19"
Set-TextValue "D4"  "5.748"
Set-TextValue "D5"  "0.05767"
Set-TextValue "D6"  "3.414"
Set-TextValue "D7"  "6.470"
Set-TextValue "D8"  "1.315"
Set-TextValue "D9"  "0.8007"
Set-TextValue "D11" "0.07628"
Set-TextValue "D12" "0.03268"
Set-TextValue "D13" "0.02974"
Set-TextValue "D14" "0.09245"
Set-TextValue "D15" "0.001668"
Set-TextValue "D16" "3.266"
Set-TextValue "D17" "0.04740"
Set-TextValue "D18" "0.0005995"
Set-TextValue "D19" "0.006226"
Set-TextValue "D20" "0.005370"
Set-TextValue "D25" "0.3323"
Set-TextValue "D26" "0.1277"
Set-TextValue "D27" "0.0006733"
Set-TextValue "D40" "0.04271"
Set-TextValue "D41" "0.007140"

# Rows 42/43: the two coins swapped rank order (CEJI <-> BKEXToken) and
# got refreshed prices.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1053"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003332"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.009540"
Set-TextValue "D46" "0.00005633"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "D48" "0.7858"
Set-TextValue "D49" "0.09705"
Set-TextValue "D50" "0.00002102"
